$wb = $excel.ActiveWorkbook

# --- 1. Update the panel-query timestamps on the "data" sheet (F2:F17) ---
$dataSheet = $wb.Worksheets.Item("data")

$newTimes = @(
    "2021-10-05 14:21:22.811435",
    "2021-10-05 14:21:22.811444",
    "2021-10-05 14:21:22.811447",
    "2021-10-05 14:21:22.811450",
    "2021-10-05 14:21:22.811453",
    "2021-10-05 14:21:22.811456",
    "2021-10-05 14:21:22.811459",
    "2021-10-05 14:21:22.811462",
    "2021-10-05 14:21:22.811465",
    "2021-10-05 14:21:22.811468",
    "2021-10-05 14:21:22.811471",
    "2021-10-05 14:21:22.811474",
    "2021-10-05 14:21:22.811476",
    "2021-10-05 14:21:22.811479",
    "2021-10-05 14:21:22.811482",
    "2021-10-05 14:21:22.811485"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = 2 + $i
    $dataSheet.Range("F$row").Value = $newTimes[$i]
}

# --- 2. Add the new "metadata" sheet after "data" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$metaSheet = $wb.Worksheets.Add($null, $lastSheet)
$metaSheet.Name = "metadata"

# Header row (bold, thin border, centered h/v) - mirrors the "data" sheet header style
$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row 2
$a2 = $metaSheet.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Value = 0

$metaSheet.Range("B2").Value = "Left Ventricular Noncompaction Cardiomyopathy"
$metaSheet.Range("C2").Value = 238
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.4"
$metaSheet.Range("E2").Value = "2020-12-02T16:46:32.432500Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:22.807832"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/238/?format=json"

# Restore the active sheet/selection to the "data" sheet, matching the original file
$dataSheet.Activate()
$dataSheet.Range("A1").Select()
